# Weekly refresh of the price-report data block (rows 12-139, columns A:R).
# Each "week" is represented by a pair of rows (Primera / Segunda quality).
# The whole block of weeks is shifted down by one slot: the data that used
# to live two rows below now lives here, a brand-new week is inserted at
# the very top (rows 12-13), and the week that used to be last (rows
# 138-139) is appended as two new rows (140-141) at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the two rows that will fall off the bottom of the shifted block,
# and the date-cell number format, before anything gets overwritten.
$lastTwoRows = $ws.Range("A138:R139").Value2
$dateFormat = $ws.Cells.Item(138, 4).NumberFormat

# Shift every week down by one slot: new row r (14..139) = old row (r-2).
$ws.Range("A14:R139").Value2 = $ws.Range("A12:R137").Value2

# The newly freed first slot (rows 12-13) gets a brand-new report date.
$ws.Cells.Item(12, 4).Value2 = 44490
$ws.Cells.Item(13, 4).Value2 = 44490

# Re-append the week that fell off the end as two brand-new rows.
$ws.Range("A140:R141").Value2 = $lastTwoRows
$ws.Cells.Item(140, 4).NumberFormat = $dateFormat
$ws.Cells.Item(141, 4).NumberFormat = $dateFormat
